$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 195 (ticker "KODK"), shifting all rows below it up by one.
$ws.Rows.Item(195).Delete()
